$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text runs: only the specific run's text changes) ---
# A8: "Volume 31   Number  49" -> "...50"  (the run containing "49")
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "50"

# C9: "Report Covering the Week  12/2/2024  Through  12/8/2024"
#     -> "...12/9/2024  Through  12/15/2024"
$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "12/9/2024"
$c9.Characters(47, 9).Text = "12/15/2024"

# --- Crime-stat table updates (rows 15-30) ---

# Row 22: C/D/E/F/G/H change; D22 and E22 switch from numbers to the
# "0" / "***.*" placeholder text (the same shared strings used by
# D15/E15). Force text entry, then copy D15/E15's format (style 13)
# onto D22/E22 so the saved style matches the rest of the sheet.
$ws.Range("C22").Value = 1

$d22 = $ws.Range("D22")
$d22.NumberFormat = "@"
$d22.Value = "0"
$ws.Range("D15").Copy()
$d22.PasteSpecial(-4122)

$e22 = $ws.Range("E22")
$e22.NumberFormat = "@"
$e22.Value = "***.*"
$ws.Range("E15").Copy()
$e22.PasteSpecial(-4122)

$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 0
$ws.Range("L22").Value = 7.407407407407
$ws.Range("M22").Value = 38.095238095238

# Row 27: D27 and E27 switch from the "0"/"***.*" placeholder text back
# to real numbers, taking on the normal numeric style of their column
# (copied from D18/E18, which keep that same style throughout this edit).
$ws.Range("C27").Value = 1

$ws.Range("D18").Copy()
$d27 = $ws.Range("D27")
$d27.PasteSpecial(-4122)
$d27.Value = 1

$ws.Range("E18").Copy()
$e27 = $ws.Range("E27")
$e27.PasteSpecial(-4122)
$e27.Value = 0

# --- Remaining plain numeric updates ---
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 6
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 44
$ws.Range("K15").Value = 10
$ws.Range("L15").Value = 51.724137931034
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -44.303797468354
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 16
$ws.Range("E16").Value = -68.75
$ws.Range("F16").Value = 34
$ws.Range("G16").Value = 57
$ws.Range("H16").Value = -40.350877192982
$ws.Range("I16").Value = 603
$ws.Range("J16").Value = 642
$ws.Range("K16").Value = -6.074766355140
$ws.Range("L16").Value = 12.290502793296
$ws.Range("M16").Value = 39.583333333333
$ws.Range("N16").Value = -67.138964577656
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 29
$ws.Range("E17").Value = -55.172413793103
$ws.Range("F17").Value = 73
$ws.Range("G17").Value = 81
$ws.Range("H17").Value = -9.876543209876
$ws.Range("I17").Value = 983
$ws.Range("J17").Value = 986
$ws.Range("K17").Value = -0.304259634888
$ws.Range("L17").Value = 19.732034104750
$ws.Range("M17").Value = 136.867469879518
$ws.Range("N17").Value = -7.088846880907
$ws.Range("C18").Value = 7
$ws.Range("E18").Value = 16.666666666666
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -7.407407407407
$ws.Range("I18").Value = 341
$ws.Range("J18").Value = 306
$ws.Range("K18").Value = 11.437908496732
$ws.Range("L18").Value = -2.571428571428
$ws.Range("M18").Value = 75.773195876288
$ws.Range("N18").Value = -73.829623944742
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -7.142857142857
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -6.779661016949
$ws.Range("I19").Value = 926
$ws.Range("J19").Value = 716
$ws.Range("K19").Value = 29.329608938547
$ws.Range("L19").Value = 29.691876750700
$ws.Range("M19").Value = 121.002386634845
$ws.Range("N19").Value = 25.474254742547
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -30
$ws.Range("I20").Value = 226
$ws.Range("J20").Value = 326
$ws.Range("K20").Value = -30.674846625766
$ws.Range("L20").Value = -22.866894197952
$ws.Range("M20").Value = 89.915966386554
$ws.Range("N20").Value = -64.353312302839
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 71
$ws.Range("E21").Value = -39.436619718309
$ws.Range("F21").Value = 207
$ws.Range("G21").Value = 250
$ws.Range("H21").Value = -17.2
$ws.Range("I21").Value = 3137
$ws.Range("J21").Value = 3028
$ws.Range("K21").Value = 3.599735799207
$ws.Range("L21").Value = 13.700616165277
$ws.Range("M21").Value = 94.241486068111
$ws.Range("N21").Value = -45.099754987749
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 0
$ws.Range("L22").Value = 7.407407407407
$ws.Range("M22").Value = 38.095238095238
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 13
$ws.Range("E23").Value = -61.538461538461
$ws.Range("F23").Value = 34
$ws.Range("G23").Value = 39
$ws.Range("H23").Value = -12.820512820512
$ws.Range("I23").Value = 466
$ws.Range("J23").Value = 470
$ws.Range("K23").Value = -0.851063829787
$ws.Range("L23").Value = 21.989528795811
$ws.Range("M23").Value = 62.937062937062
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -12.5
$ws.Range("F24").Value = 134
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = 16.521739130434
$ws.Range("I24").Value = 1660
$ws.Range("J24").Value = 1558
$ws.Range("K24").Value = 6.546854942233
$ws.Range("L24").Value = -10.656620021528
$ws.Range("M24").Value = 21.079504011670
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = -23.529411764705
$ws.Range("F25").Value = 62
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = 12.727272727272
$ws.Range("I25").Value = 732
$ws.Range("J25").Value = 611
$ws.Range("K25").Value = 19.803600654664
$ws.Range("L25").Value = -20.864864864864
$ws.Range("C26").Value = 24
$ws.Range("D26").Value = 22
$ws.Range("E26").Value = 9.090909090909
$ws.Range("F26").Value = 82
$ws.Range("H26").Value = -2.380952380952
$ws.Range("I26").Value = 1210
$ws.Range("J26").Value = 1087
$ws.Range("K26").Value = 11.315547378104
$ws.Range("L26").Value = 20.879120879120
$ws.Range("M26").Value = 10.401459854014
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 75
$ws.Range("I27").Value = 57
$ws.Range("J27").Value = 59
$ws.Range("K27").Value = -3.389830508474
$ws.Range("L27").Value = 18.75
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 250
$ws.Range("I28").Value = 154
$ws.Range("J28").Value = 104
$ws.Range("K28").Value = 48.076923076923
$ws.Range("L28").Value = 81.176470588235
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 200
$ws.Range("I29").Value = 51
$ws.Range("K29").Value = 50
$ws.Range("L29").Value = -20.3125
$ws.Range("M29").Value = -13.559322033898
$ws.Range("N29").Value = -75.714285714285
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 200
$ws.Range("I30").Value = 46
$ws.Range("K30").Value = 43.75
$ws.Range("L30").Value = -9.803921568627
$ws.Range("M30").Value = -9.803921568627
$ws.Range("N30").Value = -75.661375661375
